# Fix mojibake "Â±" -> "±" in the f1_score_weighted / training_time / test_time columns
# (B2:D17) for the automl sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "0.552 (0.551 ± 0.001)"
$ws.Range("C2").Value = "00:01:36 (00:02:40 ± 00:01:26)"
$ws.Range("D2").Value = "00:00:01 (00:00:07 ± 00:00:02)"
$ws.Range("B3").Value = "0.615 (0.552 ± 0.024)"
$ws.Range("C3").Value = "00:00:14 (00:00:17 ± 00:00:01)"
$ws.Range("D3").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B4").Value = "0.535 (0.466 ± 0.037)"
$ws.Range("C4").Value = "00:00:27 (00:00:36 ± 00:00:12)"
$ws.Range("D4").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B5").Value = "0.589 (0.529 ± 0.029)"
$ws.Range("C5").Value = "00:05:06 (00:05:13 ± 00:00:04)"
$ws.Range("D5").Value = "00:00:01 (00:00:02 ± 00:00:01)"
$ws.Range("B6").Value = "0.631 (0.569 ± 0.031)"
$ws.Range("C6").Value = "00:04:56 (00:05:00 ± 00:00:02)"
$ws.Range("D6").Value = "00:00:00 (00:00:01 ± 00:00:00)"
$ws.Range("B7").Value = "0.611 (0.564 ± 0.026)"
$ws.Range("C7").Value = "00:05:01 (00:05:04 ± 00:00:02)"
$ws.Range("D7").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B8").Value = "0.613 (0.559 ± 0.033)"
$ws.Range("C8").Value = "00:04:59 (00:06:27 ± 00:02:14)"
$ws.Range("D8").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B9").Value = "0.604 (0.552 ± 0.027)"
$ws.Range("C9").Value = "00:04:59 (00:05:00 ± 00:00:00)"
$ws.Range("D9").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B10").Value = "0.618 (0.568 ± 0.029)"
$ws.Range("C10").Value = "00:04:29 (00:04:29 ± 00:00:00)"
$ws.Range("D10").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B11").Value = "0.521 (0.348 ± 0.095)"
$ws.Range("C11").Value = "00:05:00 (00:05:05 ± 00:00:01)"
$ws.Range("D11").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B12").Value = "0.461 (0.415 ± 0.029)"
$ws.Range("C12").Value = "00:02:15 (00:03:00 ± 00:00:55)"
$ws.Range("D12").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B13").Value = "0.591 (0.550 ± 0.026)"
$ws.Range("C13").Value = "00:00:01 (00:00:01 ± 00:00:00)"
$ws.Range("D13").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B14").Value = "0.618 (0.565 ± 0.029)"
$ws.Range("C14").Value = "00:00:23 (00:00:25 ± 00:00:01)"
$ws.Range("D14").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B15").Value = "0.595 (0.564 ± 0.037)"
$ws.Range("C15").Value = "00:02:13 (00:04:13 ± 00:01:09)"
$ws.Range("D15").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B16").Value = "0.619 (0.562 ± 0.035)"
$ws.Range("C16").Value = "00:00:07 (00:00:08 ± 00:00:00)"
$ws.Range("D16").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B17").Value = "0.613 (0.561 ± 0.028)"
$ws.Range("C17").Value = "00:05:00 (00:05:14 ± 00:00:17)"
$ws.Range("D17").Value = "00:00:00 (00:00:00 ± 00:00:00)"
